# "Generate Report for Handback" - refresh the handback/handoff timestamps
# that get written each time the localization report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for the d125a089... file (row 2)
$wsOverview.Range("G2").Value = "2016-11-14 06:53:57"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for
# the d125a089... file (row 2)
$wsZhCn.Range("H2").Value = "2016-11-14 06:53:44"
$wsZhCn.Range("K2").Value = "2016-11-14 06:54:33"

# de-de: same two columns for the d125a089... file (row 2). H2 mirrors the
# Overview sheet's generate date, K2 is the handback datetime.
$wsDeDe.Range("H2").Value = "2016-11-14 06:53:57"
$wsDeDe.Range("K2").Value = "2016-11-14 06:54:52"
